$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.311.56"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.588.96"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.42"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.812.37"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "1.587.63"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "26.317.65"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "211.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "1.313.38"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.10"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.34%  "
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.769"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "1.724.97"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.48"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0981"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.47%  "
